# The commit swaps the presentation's theme palette: the deck's live theme
# (ppt/theme/theme2.xml, the "Integral"/Red Violet scheme referenced by the
# slide master) is replaced with the stock "Office Theme" palette (which
# used to live in ppt/theme/theme1.xml, used only by the notes master).
#
# PowerPoint's object model doesn't give VBA/COM a way to rename a theme or
# swap the raw part that backs it, but it does let you repaint every slot of
# the live color scheme -- which is the entire substance of this diff, since
# the font scheme and format scheme are byte-for-byte identical between the
# two themes and only the twelve clrScheme colors (plus the theme's display
# name, which isn't settable from the object model either) actually differ.
# Re-coloring the scheme in place reproduces that visible change exactly.

function BGR($r, $g, $b) {
    # COM RGB values come back/go in as 0x00BBGGRR (blue in the high byte),
    # not 0x00RRGGBB, so pack the bytes in that order.
    return ($b * 65536) + ($g * 256) + $r
}

$p  = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

# Target palette = the stock Office Theme scheme (dk1/lt1/dk2/lt2/
# accent1-6/hlink/folHlink, in that COM Colors() index order).
$cs.Colors(1).RGB  = (BGR 0x00 0x00 0x00)   # dk1
$cs.Colors(2).RGB  = (BGR 0xFF 0xFF 0xFF)   # lt1
$cs.Colors(3).RGB  = (BGR 0x44 0x54 0x6A)   # dk2
$cs.Colors(4).RGB  = (BGR 0xE7 0xE6 0xE6)   # lt2
$cs.Colors(5).RGB  = (BGR 0x5B 0x9B 0xD5)   # accent1
$cs.Colors(6).RGB  = (BGR 0xED 0x7D 0x31)   # accent2
$cs.Colors(7).RGB  = (BGR 0xA5 0xA5 0xA5)   # accent3
$cs.Colors(8).RGB  = (BGR 0xFF 0xC0 0x00)   # accent4
$cs.Colors(9).RGB  = (BGR 0x44 0x72 0xC4)   # accent5
$cs.Colors(10).RGB = (BGR 0x70 0xAD 0x47)   # accent6
$cs.Colors(11).RGB = (BGR 0x05 0x63 0xC1)   # hlink
$cs.Colors(12).RGB = (BGR 0x95 0x4F 0x72)   # folHlink
